$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.919643
$ws.Range("H2").Value = 53.75892899999999
$ws.Range("I2").Value = 0.8982899767221961
$ws.Range("J2").Value = 0.8982899767221962
$ws.Range("M2").Value = 1.334383666666667
$ws.Range("N2").Value = 4.003151
$ws.Range("O2").Value = 0.1312069045987744
$ws.Range("P2").Value = 0.1312069045987744
$ws.Range("Q2").Value = 23.91167893169766
$ws.Range("R2").Value = 215.205110385279
$ws.Range("S2").Value = 0.1178618472778244
$ws.Range("T2").Value = 0.1178618472778245
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.919643
$ws.Range("H3").Value = 53.75892899999999
$ws.Range("I3").Value = 0.8982899767221961
$ws.Range("J3").Value = 0.8982899767221962
$ws.Range("N3").Value = 7.432386999999999
$ws.Range("O3").Value = 0.2436032245723858
$ws.Range("P3").Value = 0.2436032245723858
$ws.Range("Q3").Value = 44.39524055928032
$ws.Range("R3").Value = 399.5571650335229
$ws.Range("S3").Value = 0.2188263349305803
$ws.Range("T3").Value = 0.2188263349305804
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.919643
$ws.Range("H4").Value = 53.75892899999999
$ws.Range("I4").Value = 0.8982899767221961
$ws.Range("J4").Value = 0.8982899767221962
$ws.Range("M4").Value = 6.358226000000001
$ws.Range("N4").Value = 19.074678
$ws.Range("O4").Value = 0.6251898708288398
$ws.Range("P4").Value = 0.6251898708288398
$ws.Range("Q4").Value = 113.937140033318
$ws.Range("R4").Value = 1025.434260299862
$ws.Range("S4").Value = 0.5616017945137913
$ws.Range("T4").Value = 0.5616017945137914
$ws.Range("I5").Value = 0.06812533974785755
$ws.Range("J5").Value = 0.06812533974785755
$ws.Range("M5").Value = 1.334383666666667
$ws.Range("N5").Value = 4.003151
$ws.Range("O5").Value = 0.1312069045987744
$ws.Range("P5").Value = 0.1312069045987744
$ws.Range("Q5").Value = 1.813435854096555
$ws.Range("R5").Value = 16.320922686869
$ws.Range("S5").Value = 0.008938514953056238
$ws.Range("T5").Value = 0.008938514953056238
$ws.Range("I6").Value = 0.06812533974785755
$ws.Range("J6").Value = 0.06812533974785755
$ws.Range("N6").Value = 7.432386999999999
$ws.Range("O6").Value = 0.2436032245723858
$ws.Range("P6").Value = 0.2436032245723858
$ws.Range("Q6").Value = 3.366887001594777
$ws.Range("S6").Value = 0.01659555243766742
$ws.Range("T6").Value = 0.01659555243766743
$ws.Range("I7").Value = 0.06812533974785755
$ws.Range("J7").Value = 0.06812533974785755
$ws.Range("M7").Value = 6.358226000000001
$ws.Range("N7").Value = 19.074678
$ws.Range("O7").Value = 0.6251898708288398
$ws.Range("P7").Value = 0.6251898708288398
$ws.Range("Q7").Value = 8.640869402764668
$ws.Range("R7").Value = 77.76782462488201
$ws.Range("S7").Value = 0.04259127235713388
$ws.Range("T7").Value = 0.04259127235713388
$ws.Range("G8").Value = 0.669968
$ws.Range("H8").Value = 2.009904
$ws.Range("I8").Value = 0.03358468352994624
$ws.Range("J8").Value = 0.03358468352994624
$ws.Range("M8").Value = 1.334383666666667
$ws.Range("N8").Value = 4.003151
$ws.Range("O8").Value = 0.1312069045987744
$ws.Range("P8").Value = 0.1312069045987744
$ws.Range("Q8").Value = 0.8939943563893333
$ws.Range("R8").Value = 8.045949207504
$ws.Range("S8").Value = 0.004406542367893685
$ws.Range("T8").Value = 0.004406542367893686
$ws.Range("G9").Value = 0.669968
$ws.Range("H9").Value = 2.009904
$ws.Range("I9").Value = 0.03358468352994624
$ws.Range("J9").Value = 0.03358468352994624
$ws.Range("N9").Value = 7.432386999999999
$ws.Range("O9").Value = 0.2436032245723858
$ws.Range("P9").Value = 0.2436032245723858
$ws.Range("Q9").Value = 1.659820484538666
$ws.Range("R9").Value = 14.938384360848
$ws.Range("S9").Value = 0.008181337204137999
$ws.Range("T9").Value = 0.008181337204138002
$ws.Range("G10").Value = 0.669968
$ws.Range("H10").Value = 2.009904
$ws.Range("I10").Value = 0.03358468352994624
$ws.Range("J10").Value = 0.03358468352994624
$ws.Range("M10").Value = 6.358226000000001
$ws.Range("N10").Value = 19.074678
$ws.Range("O10").Value = 0.6251898708288398
$ws.Range("P10").Value = 0.6251898708288398
$ws.Range("Q10").Value = 4.259807956768001
$ws.Range("R10").Value = 38.33827161091201
$ws.Range("S10").Value = 0.02099680395791455
$ws.Range("T10").Value = 0.02099680395791455